# Add the missing "LH_INT1_SAB" course option to the "choices" sheet.
#
# In the source survey/choices spreadsheet, the "choices" sheet lists the
# available courses. "LH_INT1_SAB" (Lindy Hop - Intermedios 1 - Sabados Sc)
# was missing and needs to be inserted right after the
# "LH_INT1_VIE_IBERA" row (row 18), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# Make "choices" the active sheet/tab, matching the final state of the file.
$ws.Activate() | Out-Null

# Insert a new row 19 (shifts rows 19.. down to 20..) and populate it.
$ws.Rows("19:19").Insert()
$ws.Range("A19").Value = "course"
$ws.Range("B19").Value = "LH_INT1_SAB"
$ws.Range("C19").Value = "Lindy Hop - Intermedios 1 - Sábados Sc"

# Restore the expected selection/page setup state on the sheet.
$ws.Range("C17").Select() | Out-Null
$ws.PageSetup.Orientation = 1
